$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a literal TEXT value into a cell without Excel re-interpreting
# numeric-looking strings (e.g. "6.90", "0.999") as numbers. We build the text
# via a formula (so the computed type is Text), copy it, then paste-special
# only the Values into the destination cell - this preserves the exact string.
$scratch = $ws.Range("Z1")
function Set-TextValue {
    param([string]$cellRef, [string]$text)
    $escaped = $text.Replace('"', '""')
    $scratch.Formula = '="' + $escaped + '"'
    $scratch.Copy() | Out-Null
    $ws.Range($cellRef).PasteSpecial(-4163) | Out-Null
}

Set-TextValue "D2" "67.048.42"
Set-TextValue "E2" "  -3.88%  "
Set-TextValue "D3" "3.524.82"
Set-TextValue "E3" "  -4.38%  "
Set-TextValue "D4" "0.999"
Set-TextValue "E4" "  -0.08%  "
Set-TextValue "D5" "608.65"
Set-TextValue "E5" "  -6.01%  "
Set-TextValue "D6" "152.04"
Set-TextValue "E6" "  -5.55%  "
Set-TextValue "D7" "3.525.39"
Set-TextValue "E7" "  -4.33%  "
Set-TextValue "E8" "  +0.00%  "
Set-TextValue "D9" "0.482"
Set-TextValue "E9" "  -4.16%  "
Set-TextValue "E10" "  -4.54%  "
Set-TextValue "D11" "6.90"
Set-TextValue "E11" "  -3.73%  "
Set-TextValue "E12" "  -4.00%  "
Set-TextValue "D13" "0.0000220"
Set-TextValue "E13" "  -5.13%  "
Set-TextValue "D14" "4.120.24"
Set-TextValue "E14" "  -4.30%  "
Set-TextValue "D15" "31.62"
Set-TextValue "E15" "  -3.52%  "
Set-TextValue "D16" "3.525.69"
Set-TextValue "E16" "  -5.02%  "
Set-TextValue "D17" "66.952.53"
Set-TextValue "E17" "  -3.99%  "
Set-TextValue "E18" "  +0.21%  "
Set-TextValue "D19" "6.30"
Set-TextValue "E19" "  -3.43%  "
Set-TextValue "D20" "15.31"
Set-TextValue "E20" "  -4.89%  "
Set-TextValue "D21" "444.87"
Set-TextValue "E21" "  -5.36%  "
Set-TextValue "D22" "9.05"
Set-TextValue "E22" "  -13.01%  "
Set-TextValue "D23" "0.630"
Set-TextValue "E23" "  -3.28%  "
Set-TextValue "D24" "77.53"
Set-TextValue "E24" "  -2.78%  "
Set-TextValue "E25" "  +0.03%  "
Set-TextValue "D26" "3.660.56"
Set-TextValue "E26" "  -4.44%  "
Set-TextValue "D27" "0.0000122"
Set-TextValue "E27" "  -2.86%  "
Set-TextValue "D28" "10.17"
Set-TextValue "E28" "  -7.59%  "
Set-TextValue "D29" "8.18"
Set-TextValue "E29" "  -10.56%  "
Set-TextValue "D30" "2.52"
Set-TextValue "E30" "  -4.80%  "
Set-TextValue "D31" "1.62"
Set-TextValue "E31" "  -6.12%  "
Set-TextValue "D32" "0.998"
Set-TextValue "E33" "  -1.32%  "
Set-TextValue "D34" "25.65"
Set-TextValue "E34" "  -4.48%  "
Set-TextValue "D35" "6.16"
Set-TextValue "E35" "  -6.29%  "
Set-TextValue "D36" "1.87"
Set-TextValue "E36" "  -7.09%  "
Set-TextValue "D37" "3.512.58"
Set-TextValue "E37" "  -4.61%  "
Set-TextValue "D38" "8.02"
Set-TextValue "E38" "  -5.37%  "
Set-TextValue "E39" "  +0.11%  "
Set-TextValue "D40" "0.998"
Set-TextValue "E40" "  -0.18%  "
Set-TextValue "D41" "172.84"
Set-TextValue "E41" "  -3.46%  "
Set-TextValue "D42" "2.14"
Set-TextValue "E42" "  -4.45%  "
Set-TextValue "D43" "5.55"
Set-TextValue "E43" "  -6.05%  "
Set-TextValue "D44" "0.0862"
Set-TextValue "E44" "  -4.60%  "
Set-TextValue "D45" "0.891"
Set-TextValue "E45" "  -3.93%  "
Set-TextValue "D46" "45.23"
Set-TextValue "E46" "  -4.26%  "
Set-TextValue "D47" "27.21"
Set-TextValue "E47" "  -6.72%  "
Set-TextValue "E48" "  -6.77%  "
Set-TextValue "E49" "  -1.68%  "
Set-TextValue "D50" "7.57"
Set-TextValue "E50" "  -3.70%  "
Set-TextValue "E51" "  -5.71%  "

$scratch.ClearContents() | Out-Null
$excel.CutCopyMode = $false
